$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 ---
$ws.Range("D9").Value = "x"
$ws.Range("H9").Value = "24/8/2022"
$ws.Range("I9").Value = "25/08/2022"

# --- Row 10 ---
$ws.Range("D10").Value = "x"
$ws.Range("H10").Value = "25/08/2022"
$ws.Range("I10").Value = "'02/09/2022"

# --- Row 11 ---
$ws.Range("E11").Value = "x"
$ws.Range("H11").Value = "'03/09/2022"
$ws.Range("I11").Value = "'05/09/2022"

# --- Row 12 ---
$ws.Range("D12").Value = "x"
$ws.Range("E12").Value = "x"
$ws.Range("H12").Value = "'06/09/2022"
$ws.Range("I12").Value = "'10/09/2022"

# --- Row 13 ---
$ws.Range("D13").Value = "x"
$ws.Range("E13").Value = "x"
$ws.Range("H13").Value = "'11/09/2022"
$ws.Range("I13").Value = "'18/09/2022"

# --- Row 14 ---
$ws.Range("H14").Value = "'18/09/2022"

# --- View/selection state: scroll + active cell ---
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 4
    $win.ScrollColumn = 1
} catch {
}
$ws.Range("H28").Select() | Out-Null
